# -----------------------------------------------------------------------
# 1. Insert a new "Player Info" worksheet as the very first sheet, with
#    player metadata for player 4716.
# -----------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

$info = $wb.Worksheets.Add($battingSheet)
$info.Name = "Player Info"

# the player ID looks numeric ("4716") - force it to be stored as text,
# matching the rest of the workbook which stores every value as text
$info.Range("A2").NumberFormat = "@"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$info.Range("A2").Value = "4716"
$info.Range("B2").Value = "Raymon Anton Reifer"
$info.Range("C2").Value = "Left Handed"
$info.Range("D2").Value = "Left Arm Medium Fast"

# Header styling to match the header row styling used on the other sheets
# (bold font, thin box border, centered/top aligned)
$hdr = $info.Range("A1:D1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108   # xlCenter
$hdr.VerticalAlignment = -4160     # xlTop
$hdr.Borders.LineStyle = 1
$hdr.Borders.Weight = 2

$info.Range("A1").Select()

# -----------------------------------------------------------------------
# 2. ODI Batting sheet: rename MATCH_CARD_LINK -> MATCH_CODE and replace
#    the full scorecard URL values with just the bare match code.
# -----------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D2:D6").NumberFormat = "@"   # keep the match codes as text
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").Value = "4293"
$batting.Range("D3").Value = "4296"
$batting.Range("D4").Value = "4443"
$batting.Range("D5").Value = "4445"
$batting.Range("D6").Value = "4447"

# -----------------------------------------------------------------------
# 3. ODI Bowling sheet: rename MATCH_CARD_LINK -> MATCH_CODE and replace
#    the full scorecard URL values with just the bare match code.
# -----------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B2:B5").NumberFormat = "@"   # keep the match codes as text
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").Value = "4293"
$bowling.Range("B3").Value = "4296"
$bowling.Range("B4").Value = "4445"
$bowling.Range("B5").Value = "4447"

Write-Host "Workbook updated: sheets =" ($wb.Worksheets | ForEach-Object { $_.Name }) -join ", "
